$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 89.38217433333334
$ws.Range("H2").Value = 268.146523
$ws.Range("I2").Value = 0.2143552015363441
$ws.Range("J2").Value = 0.2175965347165783
$ws.Range("M2").Value = 4.631270333333333
$ws.Range("N2").Value = 13.893811
$ws.Range("O2").Value = 0.05846361049715151
$ws.Range("P2").Value = 0.0596002562356855
$ws.Range("Q2").Value = 413.9530123187948
$ws.Range("R2").Value = 3725.577110869153
$ws.Range("S2").Value = 0.01253197901065923
$ws.Range("T2").Value = 0.0129688092251053
$ws.Range("G3").Value = 89.38217433333334
$ws.Range("H3").Value = 268.146523
$ws.Range("I3").Value = 0.2143552015363441
$ws.Range("J3").Value = 0.2175965347165783
$ws.Range("O3").Value = 0.01766942965546306
$ws.Range("P3").Value = 0.01801295756537869
$ws.Range("Q3").Value = 125.1088252955289
$ws.Range("R3").Value = 1125.97942765976
$ws.Range("S3").Value = 0.00378753415482904
$ws.Range("T3").Value = 0.003919557146223177
$ws.Range("G4").Value = 89.38217433333334
$ws.Range("H4").Value = 268.146523
$ws.Range("I4").Value = 0.2143552015363441
$ws.Range("J4").Value = 0.2175965347165783
$ws.Range("M4").Value = 33.32967466666667
$ws.Range("N4").Value = 99.989024
$ws.Range("O4").Value = 0.4207426855832669
$ws.Range("P4").Value = 0.428922737696382
$ws.Range("Q4").Value = 2979.078791529284
$ws.Range("R4").Value = 26811.70912376355
$ws.Range("S4").Value = 0.09018838316314383
$ws.Range("T4").Value = 0.09333210138388061
$ws.Range("G5").Value = 89.38217433333334
$ws.Range("H5").Value = 268.146523
$ws.Range("I5").Value = 0.2143552015363441
$ws.Range("J5").Value = 0.2175965347165783
$ws.Range("M5").Value = 4.5322385
$ws.Range("N5").Value = 9.064477
$ws.Range("O5").Value = 0.05721346569581108
$ws.Range("P5").Value = 0.03888387079991788
$ws.Range("Q5").Value = 405.1013317272452
$ws.Range("R5").Value = 2430.607990363471
$ws.Range("S5").Value = 0.01226400396981829
$ws.Range("T5").Value = 0.008460995542429276
$ws.Range("G6").Value = 89.38217433333334
$ws.Range("H6").Value = 268.146523
$ws.Range("I6").Value = 0.2143552015363441
$ws.Range("J6").Value = 0.2175965347165783
$ws.Range("M6").Value = 35.3234
$ws.Range("N6").Value = 105.9702
$ws.Range("O6").Value = 0.4459108085683075
$ws.Range("P6").Value = 0.454580177702636
$ws.Range("Q6").Value = 3157.282296846067
$ws.Range("R6").Value = 28415.5406716146
$ws.Range("S6").Value = 0.09558330123789371
$ws.Range("T6").Value = 0.09891507141893999
$ws.Range("I7").Value = 0.2934277926151677
$ws.Range("J7").Value = 0.2978648075949286
$ws.Range("M7").Value = 4.631270333333333
$ws.Range("N7").Value = 13.893811
$ws.Range("O7").Value = 0.05846361049715151
$ws.Range("P7").Value = 0.0596002562356855
$ws.Range("Q7").Value = 566.6544025082065
$ws.Range("R7").Value = 5099.889622573858
$ws.Range("S7").Value = 0.01715484817649211
$ws.Range("T7").Value = 0.0177528188562509
$ws.Range("I8").Value = 0.2934277926151677
$ws.Range("J8").Value = 0.2978648075949286
$ws.Range("O8").Value = 0.01766942965546306
$ws.Range("P8").Value = 0.01801295756537869
$ws.Range("S8").Value = 0.005184701740571509
$ws.Range("T8").Value = 0.005365426139427138
$ws.Range("I9").Value = 0.2934277926151677
$ws.Range("J9").Value = 0.2978648075949286
$ws.Range("M9").Value = 33.32967466666667
$ws.Range("N9").Value = 99.989024
$ws.Range("O9").Value = 0.4207426855832669
$ws.Range("P9").Value = 0.428922737696382
$ws.Range("Q9").Value = 4078.018669758696
$ws.Range("R9").Value = 36702.16802782825
$ws.Range("S9").Value = 0.1234575974896755
$ws.Range("T9").Value = 0.1277609887370228
$ws.Range("I10").Value = 0.2934277926151677
$ws.Range("J10").Value = 0.2978648075949286
$ws.Range("M10").Value = 4.5322385
$ws.Range("N10").Value = 9.064477
$ws.Range("O10").Value = 0.05721346569581108
$ws.Range("P10").Value = 0.03888387079991788
$ws.Range("Q10").Value = 554.5374625958689
$ws.Range("R10").Value = 3327.224775575213
$ws.Range("S10").Value = 0.01678802094698546
$ws.Range("T10").Value = 0.0115821366943636
$ws.Range("I11").Value = 0.2934277926151677
$ws.Range("J11").Value = 0.2978648075949286
$ws.Range("M11").Value = 35.3234
$ws.Range("N11").Value = 105.9702
$ws.Range("O11").Value = 0.4459108085683075
$ws.Range("P11").Value = 0.454580177702636
$ws.Range("Q11").Value = 4321.958918591533
$ws.Range("R11").Value = 38897.6302673238
$ws.Range("S11").Value = 0.1308426242614431
$ws.Range("T11").Value = 0.1354034371678641
$ws.Range("G12").Value = 90.33462533333334
$ws.Range("H12").Value = 271.003876
$ws.Range("I12").Value = 0.2166393574945233
$ws.Range("J12").Value = 0.2199152301234996
$ws.Range("M12").Value = 4.631270333333333
$ws.Range("N12").Value = 13.893811
$ws.Range("O12").Value = 0.05846361049715151
$ws.Range("P12").Value = 0.0596002562356855
$ws.Range("Q12").Value = 418.3640703790484
$ws.Range("R12").Value = 3765.276633411436
$ws.Range("S12").Value = 0.01266551901491297
$ws.Range("T12").Value = 0.01310700406549032
$ws.Range("G13").Value = 90.33462533333334
$ws.Range("H13").Value = 271.003876
$ws.Range("I13").Value = 0.2166393574945233
$ws.Range("J13").Value = 0.2199152301234996
$ws.Range("O13").Value = 0.01766942965546306
$ws.Range("P13").Value = 0.01801295756537869
$ws.Range("Q13").Value = 126.4419773099022
$ws.Range("R13").Value = 1137.97779578912
$ws.Range("S13").Value = 0.003827893887854194
$ws.Range("T13").Value = 0.003961323708195089
$ws.Range("G14").Value = 90.33462533333334
$ws.Range("H14").Value = 271.003876
$ws.Range("I14").Value = 0.2166393574945233
$ws.Range("J14").Value = 0.2199152301234996
$ws.Range("M14").Value = 33.32967466666667
$ws.Range("N14").Value = 99.989024
$ws.Range("O14").Value = 0.4207426855832669
$ws.Range("P14").Value = 0.428922737696382
$ws.Range("Q14").Value = 3010.823673495225
$ws.Range("R14").Value = 27097.41306145702
$ws.Range("S14").Value = 0.09114942507527915
$ws.Range("T14").Value = 0.09432664256570131
$ws.Range("G15").Value = 90.33462533333334
$ws.Range("H15").Value = 271.003876
$ws.Range("I15").Value = 0.2166393574945233
$ws.Range("J15").Value = 0.2199152301234996
$ws.Range("M15").Value = 4.5322385
$ws.Range("N15").Value = 9.064477
$ws.Range("O15").Value = 0.05721346569581108
$ws.Range("P15").Value = 0.03888387079991788
$ws.Range("Q15").Value = 409.4180668188087
$ws.Range("R15").Value = 2456.508400912852
$ws.Range("S15").Value = 0.01239468844837546
$ws.Range("T15").Value = 0.008551155395056367
$ws.Range("G16").Value = 90.33462533333334
$ws.Range("H16").Value = 271.003876
$ws.Range("I16").Value = 0.2166393574945233
$ws.Range("J16").Value = 0.2199152301234996
$ws.Range("M16").Value = 35.3234
$ws.Range("N16").Value = 105.9702
$ws.Range("O16").Value = 0.4459108085683075
$ws.Range("P16").Value = 0.454580177702636
$ws.Range("Q16").Value = 3190.926104499467
$ws.Range("R16").Value = 28718.3349404952
$ws.Range("S16").Value = 0.0966018310681015
$ws.Range("T16").Value = 0.09996910438905655
$ws.Range("G17").Value = 18.634161
$ws.Range("H17").Value = 37.268322
$ws.Range("I17").Value = 0.0446882095496985
$ws.Range("J17").Value = 0.03024263611988591
$ws.Range("M17").Value = 4.631270333333333
$ws.Range("N17").Value = 13.893811
$ws.Range("O17").Value = 0.05846361049715151
$ws.Range("P17").Value = 0.0596002562356855
$ws.Range("Q17").Value = 86.29983702585699
$ws.Range("R17").Value = 517.799022155142
$ws.Range("S17").Value = 0.00261263407692866
$ws.Range("T17").Value = 0.001802468861987798
$ws.Range("G18").Value = 18.634161
$ws.Range("H18").Value = 37.268322
$ws.Range("I18").Value = 0.0446882095496985
$ws.Range("J18").Value = 0.03024263611988591
$ws.Range("O18").Value = 0.01766942965546306
$ws.Range("P18").Value = 0.01801295756537869
$ws.Range("Q18").Value = 26.08235937944
$ws.Range("R18").Value = 156.49415627664
$ws.Range("S18").Value = 0.0007896151750669903
$ws.Range("T18").Value = 0.0005447593210926939
$ws.Range("G19").Value = 18.634161
$ws.Range("H19").Value = 37.268322
$ws.Range("I19").Value = 0.0446882095496985
$ws.Range("J19").Value = 0.03024263611988591
$ws.Range("M19").Value = 33.32967466666667
$ws.Range("N19").Value = 99.989024
$ws.Range("O19").Value = 0.4207426855832669
$ws.Range("P19").Value = 0.428922737696382
$ws.Range("Q19").Value = 621.070523816288
$ws.Range("R19").Value = 3726.423142897728
$ws.Range("S19").Value = 0.01880223729984794
$ws.Range("T19").Value = 0.01297175427969695
$ws.Range("G20").Value = 18.634161
$ws.Range("H20").Value = 37.268322
$ws.Range("I20").Value = 0.0446882095496985
$ws.Range("J20").Value = 0.03024263611988591
$ws.Range("M20").Value = 4.5322385
$ws.Range("N20").Value = 9.064477
$ws.Range("O20").Value = 0.05721346569581108
$ws.Range("P20").Value = 0.03888387079991788
$ws.Range("Q20").Value = 84.4544618993985
$ws.Range("R20").Value = 337.817847597594
$ws.Range("S20").Value = 0.002556767344078892
$ws.Range("T20").Value = 0.001175950755534574
$ws.Range("G21").Value = 18.634161
$ws.Range("H21").Value = 37.268322
$ws.Range("I21").Value = 0.0446882095496985
$ws.Range("J21").Value = 0.03024263611988591
$ws.Range("M21").Value = 35.3234
$ws.Range("N21").Value = 105.9702
$ws.Range("O21").Value = 0.4459108085683075
$ws.Range("P21").Value = 0.454580177702636
$ws.Range("Q21").Value = 658.2219226673999
$ws.Range("R21").Value = 3949.3315360044
$ws.Range("S21").Value = 0.01992695565377602
$ws.Range("T21").Value = 0.0137477029015739
$ws.Range("G22").Value = 96.27664699999998
$ws.Range("H22").Value = 288.829941
$ws.Range("I22").Value = 0.2308894388042666
$ws.Range("J22").Value = 0.2343807914451077
$ws.Range("M22").Value = 4.631270333333333
$ws.Range("N22").Value = 13.893811
$ws.Range("O22").Value = 0.05846361049715151
$ws.Range("P22").Value = 0.0596002562356855
$ws.Range("Q22").Value = 445.8831790439056
$ws.Range("R22").Value = 4012.94861139515
$ws.Range("S22").Value = 0.01349863021815854
$ws.Range("T22").Value = 0.01396915522685118
$ws.Range("G23").Value = 96.27664699999998
$ws.Range("H23").Value = 288.829941
$ws.Range("I23").Value = 0.2308894388042666
$ws.Range("J23").Value = 0.2343807914451077
$ws.Range("O23").Value = 0.01766942965546306
$ws.Range("P23").Value = 0.01801295756537869
$ws.Range("Q23").Value = 134.7590646502133
$ws.Range("R23").Value = 1212.83158185192
$ws.Range("S23").Value = 0.004079684697141332
$ws.Range("T23").Value = 0.004221891250440598
$ws.Range("G24").Value = 96.27664699999998
$ws.Range("H24").Value = 288.829941
$ws.Range("I24").Value = 0.2308894388042666
$ws.Range("J24").Value = 0.2343807914451077
$ws.Range("M24").Value = 33.32967466666667
$ws.Range("N24").Value = 99.989024
$ws.Range("O24").Value = 0.4207426855832669
$ws.Range("P24").Value = 0.428922737696382
$ws.Range("Q24").Value = 3208.869322507509
$ws.Range("R24").Value = 28879.82390256758
$ws.Range("S24").Value = 0.09714504255532047
$ws.Range("T24").Value = 0.1005312507300803
$ws.Range("G25").Value = 96.27664699999998
$ws.Range("H25").Value = 288.829941
$ws.Range("I25").Value = 0.2308894388042666
$ws.Range("J25").Value = 0.2343807914451077
$ws.Range("M25").Value = 4.5322385
$ws.Range("N25").Value = 9.064477
$ws.Range("O25").Value = 0.05721346569581108
$ws.Range("P25").Value = 0.03888387079991788
$ws.Range("Q25").Value = 436.3487261843094
$ws.Range("R25").Value = 2618.092357105857
$ws.Range("S25").Value = 0.01320998498655298
$ws.Range("T25").Value = 0.009113632412534062
$ws.Range("G26").Value = 96.27664699999998
$ws.Range("H26").Value = 288.829941
$ws.Range("I26").Value = 0.2308894388042666
$ws.Range("J26").Value = 0.2343807914451077
$ws.Range("M26").Value = 35.3234
$ws.Range("N26").Value = 105.9702
$ws.Range("O26").Value = 0.4459108085683075
$ws.Range("P26").Value = 0.454580177702636
$ws.Range("Q26").Value = 3400.818512639799
$ws.Range("R26").Value = 30607.3666137582
$ws.Range("S26").Value = 0.1029560963470933
$ws.Range("T26").Value = 0.1065448618252015
